# This script regenerates the "K" (strikeouts) column (column G) values
# in the save_data sheet. The source data pipeline changed to read the
# "K" stat directly (instead of a former "Strike#" stat), which produced
# different raw counts for column G across all data rows. Other columns
# (TB, PC, dS0, dSF, IP, I0, IF) are unaffected by this regeneration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value
$kValues = @{
    2  = 4
    3  = 1
    4  = 6
    5  = 6
    6  = 4
    7  = 8
    8  = 5
    9  = 6
    10 = 4
    11 = 4
    12 = 10
    13 = 3
    14 = 3
    15 = 9
    16 = 2
    17 = 8
    18 = 3
    19 = 5
    20 = 4
    21 = 3
    22 = 6
    23 = 4
    24 = 10
    25 = 1
    26 = 3
    27 = 3
    28 = 6
    29 = 3
    30 = 2
    31 = 2
    33 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
